# Apply updated dSF (column F) values to Sheet1, reflecting a repull/push of
# data and mean calculation, per the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row => new value for column F (dSF)
$updates = @{
    2  = -2
    3  = -1
    5  = -5
    7  = -5
    8  = 1
    9  = -2
    10 = -2
    11 = 2
    13 = -6
    14 = -1
    16 = 2
    17 = -2
    18 = -8
    19 = 0
    20 = -4
    21 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
